$d = $word.ActiveDocument

# The paragraph containing the Word field:
#   { m:'doc.html'.fromHTMLURI() }
# is built out of fldChar begin/end + instrText runs. This change rewrites
# it into plain literal-text runs (split the same way the field's instrText
# runs were split), keeping the existing bookmark in its original place.

$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $targetParagraph = $p
        break
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not locate the paragraph containing the field."
}

# Recover the paragraph's own <w:p ...> opening tag (with its original
# w:rsidR / w:rsidRDefault / w:rsidP bookkeeping attributes) so the
# rewritten paragraph keeps looking like the rest of the document.
$openTag = '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">'
try {
    $full = $targetParagraph.Range.WordOpenXML
    if ($full -match '<w:p\s+([^>]*?)>') {
        $attrs = $matches[1]
        # Drop w14:* attributes (paraId/textId) that WordOpenXML adds but
        # that are not present on the paragraph in the stored document.
        $attrs = [System.Text.RegularExpressions.Regex]::Replace($attrs, 'w14:\S+="[^"]*"\s*', '')
        $attrs = $attrs.Trim()
        if ($attrs.Length -gt 0) {
            $openTag = '<w:p ' + $attrs + '>'
        } else {
            $openTag = '<w:p>'
        }
    }
} catch {
    # Keep the fallback $openTag defined above.
}

$body = $openTag + `
    '<w:r><w:t>{</w:t></w:r>' + `
    '<w:r><w:t>m</w:t></w:r>' + `
    '<w:r><w:t>:</w:t></w:r>' + `
    "<w:r><w:t>'</w:t></w:r>" + `
    '<w:r><w:t>doc.html</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    "<w:r><w:t>'.fromHTMLURI()</w:t></w:r>" + `
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
    '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$null = $targetParagraph.Range.InsertXML($xml)
